# Generate Report for Handoff
# Adds a new tracked file (e35be574-99a1-483f-bd6c-fda89ccae153.md) as row 3
# on all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileBase = "e35be574-99a1-483f-bd6c-fda89ccae153"
$mdName = "$fileBase.md"
$mdDisplayPath = "e2e\$fileBase.md"
$commitHash = "d6b42784d35b80d9ddd1f7c5a43e77119bf5bf14"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$mdName"

# ---------------------------------------------------------------------------
# Overview sheet (row 3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $mdUrl, "", "", $mdDisplayPath)
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-02 02:46:49"
$wsOverview.Range("G3").NumberFormat = $wsOverview.Range("G2").NumberFormat

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet (row 3)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdUrl, "", "", $mdName)
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "$fileBase.01074b8ec0c46eee698ea231e346ae5c3b667b4a.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-02 02:46:44"
$wsZhCn.Range("H3").NumberFormat = $wsZhCn.Range("H2").NumberFormat
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $wsZhCn.Range("K2").NumberFormat
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet (row 3)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdUrl, "", "", $mdName)
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "$fileBase.01074b8ec0c46eee698ea231e346ae5c3b667b4a.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-02 02:46:49"
$wsDeDe.Range("H3").NumberFormat = $wsDeDe.Range("H2").NumberFormat
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $wsDeDe.Range("K2").NumberFormat
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
